$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-pulled "dSF" (column F) data - updating specific rows with refreshed values
$ws.Range("F2").Value  = -4
$ws.Range("F3").Value  = -4
$ws.Range("F7").Value  = -4
$ws.Range("F8").Value  = -3
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = -4
$ws.Range("F16").Value = -1
$ws.Range("F19").Value = -6
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = 4
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = -1
$ws.Range("F31").Value = 3
$ws.Range("F32").Value = -5
$ws.Range("F33").Value = -9
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = -3
$ws.Range("F37").Value = -2
$ws.Range("F38").Value = -7
